$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 114.55556
$ws.Range("I5").Value = 105.375
$ws.Range("J5").Value = 188
$ws.Range("K5").Value = 105.375
$ws.Range("L5").Value = 188
$ws.Range("M5").Value = 9.625
$ws.Range("N5").Value = -418
$ws.Range("H9").Value = 384807.7
$ws.Range("I9").Value = 219.5
$ws.Range("K9").Value = 219.5
$ws.Range("M9").Value = -50.5
$ws.Range("H17").Value = 2333.8462
$ws.Range("J17").Value = 2333.8462
$ws.Range("L17").Value = 7001.5386
$ws.Range("N17").Value = -7337.5386
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 224.06667
$ws.Range("I33").Value = 188
$ws.Range("J33").Value = 729
$ws.Range("K33").Value = 188
$ws.Range("L33").Value = 729
$ws.Range("M33").Value = 41
$ws.Range("N33").Value = -1187
$ws.Range("H34").Value = 9999
$ws.Range("I34").Value = 9999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 9999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -9796
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 9999
$ws.Range("I36").Value = 9999
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 9999
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -9284
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 6999.846
$ws.Range("I40").Value = 10571.286
$ws.Range("K40").Value = 10571.286
$ws.Range("M40").Value = -10396.286
$ws.Range("H47").Value = 9975
$ws.Range("J47").Value = 9975
$ws.Range("L47").Value = 9975
$ws.Range("N47").Value = -11919
$ws.Range("H53").Value = 441.42856
$ws.Range("J53").Value = 488
$ws.Range("L53").Value = 488
$ws.Range("N53").Value = -1762
$ws.Range("H62").Value = 3811.5
$ws.Range("I62").Value = 2775.875
$ws.Range("K62").Value = 2775.875
$ws.Range("M62").Value = -2151.875
$ws.Range("H65").Value = 3811.5
$ws.Range("I65").Value = 2775.875
$ws.Range("K65").Value = 13879.375
$ws.Range("M65").Value = -10759.375
$ws.Range("H92").Value = 3335
$ws.Range("I92").Value = 675.1
$ws.Range("J92").Value = 9984.75
$ws.Range("K92").Value = 675.1
$ws.Range("L92").Value = 9984.75
$ws.Range("M92").Value = 572.9
$ws.Range("N92").Value = -12480.75
$ws.Range("H107").Value = 296.9375
$ws.Range("I107").Value = 250.13333
$ws.Range("K107").Value = 250.13333
$ws.Range("M107").Value = 1669.86667
$ws.Range("H132").Value = 287055.6
$ws.Range("I132").Value = 1185.3549
$ws.Range("K132").Value = 3556.0647
$ws.Range("M132").Value = -1026.0647
$ws.Range("H137").Value = 3689.5806
$ws.Range("I137").Value = 1431.9412
$ws.Range("K137").Value = 4295.8236
$ws.Range("M137").Value = -1745.8236
$ws.Range("H138").Value = 3886.5518
$ws.Range("I138").Value = 2644.1538
$ws.Range("J138").Value = 4896
$ws.Range("K138").Value = 7932.4614
$ws.Range("L138").Value = 14688
$ws.Range("M138").Value = -2792.4614
$ws.Range("N138").Value = -24968
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1733
$ws.Range("I5").Value = 2000
$ws.Range("K5").Value = 2000
$ws.Range("M5").Value = -1888
$ws.Range("H26").Value = 4999.3335
$ws.Range("I26").Value = 4999
$ws.Range("K26").Value = 4999
$ws.Range("M26").Value = -4669
$ws.Range("H29").Value = 2380
$ws.Range("J29").Value = 2840
$ws.Range("L29").Value = 2840
$ws.Range("N29").Value = -3456
$ws.Range("H32").Value = 17378.701
$ws.Range("I32").Value = 2943.2856
$ws.Range("J32").Value = 59482
$ws.Range("K32").Value = 2943.2856
$ws.Range("L32").Value = 59482
$ws.Range("M32").Value = -2656.2856
$ws.Range("N32").Value = -60056
$ws.Range("H37").Value = 2517
$ws.Range("I37").Value = 2517
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 2517
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -2244
$ws.Range("N37").ClearContents()
$ws.Range("H42").Value = 27176.455
$ws.Range("I42").Value = 19998.75
$ws.Range("J42").Value = 31278
$ws.Range("K42").Value = 19998.75
$ws.Range("L42").Value = 31278
$ws.Range("M42").Value = -19512.75
$ws.Range("N42").Value = -32250
$ws.Range("H63").Value = 2588.7334
$ws.Range("I63").Value = 2486.0833
$ws.Range("K63").Value = 2486.0833
$ws.Range("M63").Value = -1800.0833
$ws.Range("H66").Value = 2588.7334
$ws.Range("I66").Value = 2486.0833
$ws.Range("K66").Value = 12430.4165
$ws.Range("M66").Value = -8998.416499999999
$ws.Range("H97").Value = 1041.8096
$ws.Range("I97").Value = 1041.8096
$ws.Range("K97").Value = 1041.8096
$ws.Range("M97").Value = -545.8096
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1733
$ws.Range("I4").Value = 2000
$ws.Range("K4").Value = 2000
$ws.Range("M4").Value = -1885
$ws.Range("H22").Value = 513.5625
$ws.Range("I22").Value = 513.5625
$ws.Range("K22").Value = 513.5625
$ws.Range("M22").Value = -340.5625
$ws.Range("H99").Value = 1552.4546
$ws.Range("I99").Value = 1353.8572
$ws.Range("K99").Value = 1353.8572
$ws.Range("M99").Value = 144.1428000000001
$ws.Range("H100").Value = 55749
$ws.Range("J100").Value = 62665.332
$ws.Range("L100").Value = 62665.332
$ws.Range("N100").Value = -64829.332
$ws.Range("H107").Value = 1455.6
$ws.Range("I107").Value = 1356.4286
$ws.Range("K107").Value = 1356.4286
$ws.Range("M107").Value = 563.5714
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 10009
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10009
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10009
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -10349
$ws.Range("H22").Value = 842.3
$ws.Range("I22").Value = 737.25
$ws.Range("K22").Value = 737.25
$ws.Range("M22").Value = -387.25
$ws.Range("H31").Value = 4777.932
$ws.Range("I31").Value = 3256.4595
$ws.Range("K31").Value = 3256.4595
$ws.Range("M31").Value = -2961.4595
$ws.Range("H34").Value = 4777.932
$ws.Range("I34").Value = 3256.4595
$ws.Range("K34").Value = 3256.4595
$ws.Range("M34").Value = -3054.4595
$ws.Range("H35").Value = 127229.375
$ws.Range("I35").Value = 203079.4
$ws.Range("K35").Value = 203079.4
$ws.Range("M35").Value = -202785.4
$ws.Range("H39").Value = 16732.6
$ws.Range("I39").Value = 10415.75
$ws.Range("K39").Value = 10415.75
$ws.Range("M39").Value = -10024.75
$ws.Range("H41").Value = 31686
$ws.Range("J41").Value = 37499.5
$ws.Range("L41").Value = 37499.5
$ws.Range("N41").Value = -38355.5
$ws.Range("H49").Value = 16732.6
$ws.Range("I49").Value = 10415.75
$ws.Range("K49").Value = 10415.75
$ws.Range("M49").Value = -10233.75
$ws.Range("H54").Value = 24999
$ws.Range("J54").Value = 24999
$ws.Range("L54").Value = 24999
$ws.Range("N54").Value = -26315
$ws.Range("H99").Value = 2756.9395
$ws.Range("I99").Value = 2836.0386
$ws.Range("K99").Value = 2836.0386
$ws.Range("M99").Value = -1338.0386
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350
$ws.Range("H120").Value = 43332.332
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258
$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -42620
$ws.Range("H122").Value = 2139.647
$ws.Range("J122").Value = 999
$ws.Range("L122").Value = 2997
$ws.Range("N122").Value = -7897
$ws.Range("H126").Value = 2756.9395
$ws.Range("I126").Value = 2836.0386
$ws.Range("K126").Value = 8508.1158
$ws.Range("M126").Value = -6038.1158
$ws.Range("H131").Value = 50000
$ws.Range("J131").Value = 50000
$ws.Range("L131").Value = 50000
$ws.Range("N131").Value = -60080
$ws.Range("H132").Value = 2861.0667
$ws.Range("I132").Value = 2524.4546
$ws.Range("K132").Value = 7573.3638
$ws.Range("M132").Value = -5043.3638
$ws.Range("H135").Value = 78852.5
$ws.Range("J135").Value = 78852.5
$ws.Range("L135").Value = 78852.5
$ws.Range("N135").Value = -88992.5
$ws.Range("H141").Value = 313959.88
$ws.Range("J141").Value = 313959.88
$ws.Range("L141").Value = 313959.88
$ws.Range("N141").Value = -324319.88
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 76.52631
$ws.Range("I2").Value = 54.875
$ws.Range("K2").Value = 329.25
$ws.Range("M2").Value = -216.25
$ws.Range("H23").Value = 97.55556
$ws.Range("I23").Value = 97.55556
$ws.Range("K23").Value = 292.66668
$ws.Range("M23").Value = -57.66667999999999
$ws.Range("H39").Value = 3583
$ws.Range("J39").Value = 2999.5
$ws.Range("L39").Value = 8998.5
$ws.Range("N39").Value = -9586.5
$ws.Range("H46").Value = 252.2
$ws.Range("J46").Value = 432
$ws.Range("L46").Value = 1296
$ws.Range("N46").Value = -1478
$ws.Range("H50").Value = 7988.3335
$ws.Range("I50").Value = 7559.2
$ws.Range("J50").Value = 8524.75
$ws.Range("K50").Value = 22677.6
$ws.Range("L50").Value = 25574.25
$ws.Range("M50").Value = -22196.6
$ws.Range("N50").Value = -26536.25
$ws.Range("H53").Value = 7988.3335
$ws.Range("I53").Value = 7559.2
$ws.Range("J53").Value = 8524.75
$ws.Range("K53").Value = 22677.6
$ws.Range("L53").Value = 25574.25
$ws.Range("M53").Value = -22196.6
$ws.Range("N53").Value = -26536.25
$ws.Range("H55").Value = 400213.6
$ws.Range("I55").Value = 234.5
$ws.Range("J55").Value = 666866.3
$ws.Range("K55").Value = 703.5
$ws.Range("L55").Value = 2000598.9
$ws.Range("M55").Value = -526.5
$ws.Range("N55").Value = -2000952.9
$ws.Range("H92").Value = 480.23077
$ws.Range("I92").Value = 456.25
$ws.Range("J92").Value = 518.6
$ws.Range("K92").Value = 1368.75
$ws.Range("L92").Value = 1555.8
$ws.Range("M92").Value = -120.75
$ws.Range("N92").Value = -4051.8
$ws.Range("H109").Value = 5726.7
$ws.Range("I109").Value = 2872.8333
$ws.Range("K109").Value = 8618.499899999999
$ws.Range("M109").Value = -7578.499899999999
$ws.Range("H131").Value = 67122.06
$ws.Range("I131").Value = 186123.17
$ws.Range("K131").Value = 558369.51
$ws.Range("M131").Value = -553329.51
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 935
$ws.Range("I19").Value = 935
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 935
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -647
$ws.Range("N19").ClearContents()
$ws.Range("H43").Value = 5159.8
$ws.Range("I43").Value = 5159.8
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 5159.8
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -5008.8
$ws.Range("N43").ClearContents()
$ws.Range("H59").Value = 11209.833
$ws.Range("I59").Value = 11615.777
$ws.Range("J59").Value = 9992
$ws.Range("K59").Value = 11615.777
$ws.Range("L59").Value = 9992
$ws.Range("M59").Value = -11032.777
$ws.Range("N59").Value = -11158
$ws.Range("H80").Value = 7730.5293
$ws.Range("I80").Value = 4770.1333
$ws.Range("J80").Value = 10067.685
$ws.Range("K80").Value = 4770.1333
$ws.Range("L80").Value = 10067.685
$ws.Range("M80").Value = -3772.1333
$ws.Range("N80").Value = -12063.685
$ws.Range("H83").Value = 7730.5293
$ws.Range("I83").Value = 4770.1333
$ws.Range("J83").Value = 10067.685
$ws.Range("K83").Value = 23850.6665
$ws.Range("L83").Value = 50338.425
$ws.Range("M83").Value = -18858.6665
$ws.Range("N83").Value = -60322.425
$ws.Range("H113").Value = 3456.8572
$ws.Range("I113").Value = 3456.8572
$ws.Range("K113").Value = 3456.8572
$ws.Range("M113").Value = -1286.8572
$ws.Range("H122").Value = 34578.582
$ws.Range("I122").Value = 60823.35
$ws.Range("J122").Value = 2709.9285
$ws.Range("K122").Value = 182470.05
$ws.Range("L122").Value = 8129.7855
$ws.Range("M122").Value = -180020.05
$ws.Range("N122").Value = -13029.7855
$ws.Range("H132").Value = 2156.625
$ws.Range("I132").Value = 1884.2916
$ws.Range("K132").Value = 5652.8748
$ws.Range("M132").Value = -3122.8748
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 14003
$ws.Range("J18").Value = 14003
$ws.Range("L18").Value = 14003
$ws.Range("N18").Value = -14347
$ws.Range("H20").Value = 9290
$ws.Range("H24").Value = 1100007
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 1100007
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 1100007
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -1100693
$ws.Range("H29").Value = 42999
$ws.Range("I29").Value = 39998.332
$ws.Range("J29").Value = 47500
$ws.Range("K29").Value = 39998.332
$ws.Range("L29").Value = 47500
$ws.Range("M29").Value = -39703.332
$ws.Range("N29").Value = -48090
$ws.Range("H40").Value = 4123.909
$ws.Range("I40").Value = 3656.2856
$ws.Range("K40").Value = 3656.2856
$ws.Range("M40").Value = -3520.2856
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H68").Value = 1849.9231
$ws.Range("I68").Value = 1595.3636
$ws.Range("K68").Value = 1595.3636
$ws.Range("M68").Value = -846.3635999999999
$ws.Range("H71").Value = 1849.9231
$ws.Range("I71").Value = 1595.3636
$ws.Range("K71").Value = 7976.817999999999
$ws.Range("M71").Value = -4232.817999999999
$ws.Range("H88").Value = 27363.666
$ws.Range("I88").Value = 22499.5
$ws.Range("K88").Value = 22499.5
$ws.Range("M88").Value = -22071.5
$ws.Range("H91").Value = 27363.666
$ws.Range("I91").Value = 22499.5
$ws.Range("K91").Value = 22499.5
$ws.Range("M91").Value = -21017.5
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H105").Value = 25971.666
$ws.Range("J105").Value = 25971.666
$ws.Range("L105").Value = 25971.666
$ws.Range("N105").Value = -32959.666
$ws.Range("H116").Value = 130000
$ws.Range("J116").Value = 130000
$ws.Range("L116").Value = 130000
$ws.Range("N116").Value = -139178
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H136").Value = 3924.4
$ws.Range("I136").Value = 3434
$ws.Range("J136").Value = 4660
$ws.Range("K136").Value = 10302
$ws.Range("L136").Value = 13980
$ws.Range("M136").Value = -7752
$ws.Range("N136").Value = -19080
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 899
$ws.Range("J13").Value = 899
$ws.Range("L13").Value = 899
$ws.Range("N13").Value = -1179
$ws.Range("H19").Value = 2506
$ws.Range("J19").Value = 2506
$ws.Range("L19").Value = 2506
$ws.Range("N19").Value = -2854
$ws.Range("H108").Value = 90600
$ws.Range("J108").Value = 90600
$ws.Range("L108").Value = 90600
$ws.Range("N108").Value = -98280
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -12530
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 1483.65
$ws.Range("I132").Value = 1357.5834
$ws.Range("K132").Value = 4072.7502
$ws.Range("M132").Value = -1542.7502
$ws.Range("H136").Value = 4033.5
$ws.Range("I136").Value = 4033.5
$ws.Range("K136").Value = 12100.5
$ws.Range("M136").Value = -9550.5
